$d = $word.ActiveDocument

# 1. Merge the "7" + " april" runs into a single run reading "7 april".
#    Use MatchWholeWord + a single, non-wrapping replace so the engine
#    only touches this exact occurrence (not the later "07 april" cell).
$rng1 = $d.Content
$rng1.Find.Execute("7 april", $true, $false, $false, $false, $false, `
                    $true, 0, $false, "7 april", 1) | Out-Null

# 2. Replace "Jobbat med logo" with the longer finished-logo description,
#    then add a second bold run "Fortsatt med färgpalett" right after it
#    (same run formatting: bold + size 28) inside the same table cell.
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Jobbat med logo")
if ($found2) {
    $rng2.Text = "Utvecklat en logo utifrån min tidigare design från gränssnittsdesign. Den är nu klar. "
    $rng2.Collapse(0)
    $rng2.InsertAfter("Fortsatt med färgpalett")

    # Toggle Bold off/on so the new run stays a distinct run instead of
    # silently re-merging with the identically-formatted preceding run,
    # while still ending up with the very same (bold, size 28) formatting.
    $rng2.Font.Bold = 0
    $rng2.Font.Bold = 1
    $rng2.Font.Size = 14
}
